$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 2, shifting existing rows 2-4 down to 4-6.
# Using Rows(2).Resize(2) to select a 2-row block starting at row 2, then Insert.
$ws.Rows.Item(2).Resize(2).Insert()

# Newly inserted rows inherit formatting from the row above (header row);
# clear that so the new data rows start from the default (unstyled) format,
# matching the plain cells used by the other data rows.
$ws.Range("A2:D3").ClearFormats()

# Row 2: Betonieren
$ws.Range("A2").Value = "0_5_T1_SP_GRU_EG0_3101_03_F-P-001 - Wand_Kein BA_Beton@Betonieren"
$ws.Range("B2").Value = "Beton@Betonieren"
$ws.Range("C2").Value = 45051
$ws.Range("D2").Value = 45052
$ws.Range("C2:D2").NumberFormat = $ws.Range("C4:D4").NumberFormat

# Row 3: Bewehren
$ws.Range("A3").Value = "0_5_T1_SP_GRU_EG0_3101_03_F-P-001 - Wand_Kein BA_Beton@Bewehren"
$ws.Range("B3").Value = "Beton@Bewehren"
$ws.Range("C3").Value = 45052
$ws.Range("D3").Value = 45053
$ws.Range("C3:D3").NumberFormat = $ws.Range("C4:D4").NumberFormat

$ws.Range("A1").Select()
